$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "原编号"
$ws.Range("L1").Value = "保护现状"
$ws.Range("M1").Value = "照片"

$ws.Range("K1").Select()
